$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated emmCI string values (row 2 = female, row 3 = male) ---
$ws.Range("D2").Value = "6.40 (6.35, 6.46)"
$ws.Range("D3").Value = "7.78 (7.71, 7.86)"
$ws.Range("E2").Value = "6.40 (6.35, 6.45)"
$ws.Range("E3").Value = "7.79 (7.72, 7.87)"

# --- Updated numeric statistics, row 2 ---
$ws.Range("F2").Value = 0.045806251947808324
$ws.Range("G2").Value = 0.045806251947808296
$ws.Range("H2").Value = -0.0016437045720376631
$ws.Range("I2").Value = -0.025662985979894047
$ws.Range("J2").Value = 4.0090864702492341
$ws.Range("M2").Value = 0.18090670659715877
$ws.Range("N2").Value = 0.0081154103842393053

# --- Updated numeric statistics, row 3 ---
$ws.Range("F3").Value = [double]"5.9144749045337464e-10"
$ws.Range("H3").Value = 0.0097633539683767978
$ws.Range("I3").Value = 0.12548337809205459
$ws.Range("J3").Value = 39.932137000209629
$ws.Range("M3").Value = 0.57094362994356962
$ws.Range("N3").Value = 0.075353303210206765

# --- Column J width adjustment (target stored width 10.7109375) ---
$ws.Columns.Item(10).ColumnWidth = 9.75
